$d = $word.ActiveDocument

# Locate the end of the last bullet's text ("...based on their current stash"),
# right before the trailing bookmark (_GoBack) that lives at the very end of
# the document, and collapse the range there.
$r = $d.Content
$found = $r.Find.Execute("based on their current stash", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

# Create a new paragraph right after the current last one (same list style),
# inheriting the paragraph-mark formatting of the existing bullet list.
$null = $r.InsertParagraphAfter()

# Insert the new bullet's own paragraph (style/numbering/run) as a clean XML
# fragment so the new run doesn't pick up the bold paragraph-mark formatting.
$newPara = $d.Paragraphs.Last
$r2 = $newPara.Range
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>As a user suggested, an “aggressive” Pop-Up after a user entered their recipes for a mix, to make people think a little more about steeping time</w:t></w:r></w:p>'
$null = $r2.InsertXML($frag)
